$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("G3").Value = 2.35
$ws.Range("I3").Value = 3.7
$ws.Range("J3").Value = 3.25
$ws.Range("K3").Value = 1.8
$ws.Range("M3").Value = 1.17
$ws.Range("N3").Value = 4.75
$ws.Range("O3").Value = 1.73
$ws.Range("P3").Value = 2
$ws.Range("Q3").Value = 3.5
$ws.Range("R3").Value = 1.3
$ws.Range("S3").Value = 1.78
$ws.Range("T3").Value = 2.03
$ws.Range("U3").Value = 2.63
$ws.Range("V3").Value = 1.44
$ws.Range("W3").Value = 5
$ws.Range("AC3").Value = 4.75
$ws.Range("AJ3").Value = 15
$ws.Range("AT3").Value = 2
$ws.Range("AU3").Value = 11
$ws.Range("BA3").Value = 151

# Row 4
$ws.Range("G4").Value = 1.6
$ws.Range("H4").Value = 3.7
$ws.Range("I4").Value = 6.5
$ws.Range("J4").Value = 2.2
$ws.Range("M4").Value = 1.07
$ws.Range("N4").Value = 8.5
$ws.Range("O4").Value = 1.33
$ws.Range("P4").Value = 3.25
$ws.Range("Q4").Value = 2.1
$ws.Range("R4").Value = 1.7
$ws.Range("Y4").Value = 8.5
$ws.Range("AC4").Value = 8.5
$ws.Range("AE4").Value = 21
$ws.Range("AJ4").Value = 19
$ws.Range("AP4").Value = 21
$ws.Range("AQ4").Value = 26
$ws.Range("BB4").Value = 351

# Row 5
$ws.Range("G5").Value = 4.33
$ws.Range("I5").Value = 2
$ws.Range("J5").Value = 5
$ws.Range("AK5").Value = 17
$ws.Range("AN5").Value = 6

# Row 6
$ws.Range("AG6").Value = 700

# Row 7
$ws.Range("G7").Value = 1.37
$ws.Range("L7").Value = 6
$ws.Range("U7").Value = 1.7
$ws.Range("V7").Value = 2.05
$ws.Range("X7").Value = 8
$ws.Range("Z7").Value = 9.75
$ws.Range("AF7").Value = 55
$ws.Range("AH7").Value = 25
$ws.Range("AI7").Value = 55
$ws.Range("AL7").Value = 65
$ws.Range("AM7").Value = 50
$ws.Range("AN7").Value = 3.45
$ws.Range("AO7").Value = 6.2
$ws.Range("AP7").Value = 13.5
$ws.Range("AQ7").Value = 16
$ws.Range("AR7").Value = 35
$ws.Range("AU7").Value = 7.5
$ws.Range("AY7").Value = 32
$ws.Range("AZ7").Value = 200

# Row 8
$ws.Range("G8").Value = 2.3
$ws.Range("H8").Value = 3.35
$ws.Range("I8").Value = 2.92
$ws.Range("J8").Value = 2.87
$ws.Range("L8").Value = 3.55
$ws.Range("W8").Value = 10
$ws.Range("X8").Value = 13
$ws.Range("Z8").Value = 25
$ws.Range("AA8").Value = 16.5
$ws.Range("AD8").Value = 6.6
$ws.Range("AE8").Value = 11.75
$ws.Range("AH8").Value = 10.5
$ws.Range("AI8").Value = 16
$ws.Range("AM8").Value = 27
$ws.Range("AP8").Value = 18.5
$ws.Range("AS8").Value = 200
$ws.Range("AY8").Value = 23
$ws.Range("AZ8").Value = 80
$ws.Range("BB8").Value = 300

# Row 9
$ws.Range("G9").Value = 2.12
$ws.Range("H9").Value = 3.45
$ws.Range("I9").Value = 3.15
$ws.Range("J9").Value = 2.67
$ws.Range("K9").Value = 2.15
$ws.Range("L9").Value = 3.65
$ws.Range("M9").Value = 1.05
$ws.Range("N9").Value = 8
$ws.Range("O9").Value = 1.24
$ws.Range("P9").Value = 3.65
$ws.Range("Q9").Value = 1.72
$ws.Range("R9").Value = 2.05
$ws.Range("S9").Value = 1.37
$ws.Range("T9").Value = 2.87
$ws.Range("V9").Value = 2.2
$ws.Range("W9").Value = 9.25
$ws.Range("X9").Value = 11.5
$ws.Range("Y9").Value = 8.5
$ws.Range("Z9").Value = 21
$ws.Range("AA9").Value = 15.5
$ws.Range("AC9").Value = 8
$ws.Range("AD9").Value = 6.8
$ws.Range("AH9").Value = 11
$ws.Range("AI9").Value = 17.5
$ws.Range("AJ9").Value = 11
$ws.Range("AK9").Value = 40
$ws.Range("AL9").Value = 26
$ws.Range("AM9").Value = 30
$ws.Range("AN9").Value = 4.15
$ws.Range("AO9").Value = 10.75
$ws.Range("AP9").Value = 17.5
$ws.Range("AQ9").Value = 40
$ws.Range("AR9").Value = 65
$ws.Range("AT9").Value = 2.87
$ws.Range("AU9").Value = 6.8
$ws.Range("AW9").Value = 5.2
$ws.Range("AX9").Value = 17
$ws.Range("AY9").Value = 23
$ws.Range("AZ9").Value = 80
$ws.Range("BA9").Value = 110
$ws.Range("BB9").Value = 250
